$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Issue 197 and 199 moved from OPEN to DONE
$ws.Range("G197").Value = "DONE"
$ws.Range("G199").Value = "DONE"

# New issues 200-206
$ws.Range("A200").Value = 200
$ws.Range("B200").Value = "Put in undo redo for CET editPixels"
$ws.Range("F200").Value = 43932
$ws.Range("G200").Value = "DONE"

$ws.Range("A201").Value = 201
$ws.Range("B201").Value = "Remove interface Doable and convertToRuntimeException from DialogView to somewhere more generic"
$ws.Range("F201").Value = 43932
$ws.Range("G201").Value = "OPEN"

$ws.Range("A202").Value = 202
$ws.Range("B202").Value = "there seems to be a more generic problem with undo redo for transforms … check 200 works "
$ws.Range("F202").Value = 43932
$ws.Range("G202").Value = "OPEN"

$ws.Range("A203").Value = 203
$ws.Range("B203").Value = "need to make change pixel chain thickness draggable"
$ws.Range("F203").Value = 43932
$ws.Range("G203").Value = "DONE"

$ws.Range("A204").Value = 204
$ws.Range("B204").Value = "slider bar on zoom or buttons"
$ws.Range("F204").Value = 43932
$ws.Range("G204").Value = "OPEN"

$ws.Range("A205").Value = 205
$ws.Range("B205").Value = "move control validation in EPMD to the generic validator based on the container"
$ws.Range("F205").Value = 43933
$ws.Range("G205").Value = "OPEN"

$ws.Range("A206").Value = 206
$ws.Range("B206").Value = "move show graffiti in EPMD to above show Curves"
$ws.Range("F206").Value = 43933
$ws.Range("G206").Value = "OPEN"

# Move the active selection (bottom pane) to the newly added issue row
$ws.Range("B203").Select()
